$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 163; this shifts the existing rows
# 163-167 down to 164-168 (and copies formatting, e.g. the date style
# on column D, from the row above).
$ws.Rows("163").Insert()

# Populate the newly inserted row 163 with the new record
# (Fruta, Feria Lagunitas de Puerto Montt - Durazno / Florida King / Tercera).
$ws.Range("A163").Value = 4
$ws.Range("B163").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C163").Value = "Los Lagos"
$ws.Range("D163").Value = 44505
$ws.Range("E163").Value = 10
$ws.Range("F163").Value = "Fruta"
$ws.Range("G163").Value = 100103
$ws.Range("H163").Value = "Frutos de hueso (carozo)"
$ws.Range("I163").Value = 100103004
$ws.Range("J163").Value = "Durazno"
$ws.Range("K163").Value = "Florida King"
$ws.Range("L163").Value = "Tercera"
$ws.Range("M163").Value = 100
$ws.Range("N163").Value = 18000
$ws.Range("O163").Value = 18000
$ws.Range("P163").Value = 18000
$ws.Range("Q163").Value = "$/caja 14 kilos empedrada"
$ws.Range("R163").Value = "Provincia de Limarí"
$ws.Range("S163").Value = 1286
$ws.Range("T163").Value = 14
